# Reorganización completa: limpieza de módulos antiguos, nuevas entregas y optimización

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "Datos" to "razas"
$ws.Name = "razas"

# Update header values to lowercase, snake_case style and clear the old bold/fill style
$headers = @("codigo", "nombre", "tipo_ganado", "especie", "descripcion", "comentario")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.ClearFormats()
}
